# Applies the cryptos.xlsx data refresh described in the commit diff.
# Column D ("Price") and E ("Volume(1h)") are stored as literal text in the
# sheet (t="inlineStr"), and several new prices are plain decimal numbers
# (e.g. "1.000", "20.80"). Excel normally auto-converts such text to a number
# on assignment, so NumberFormat is forced to Text ("@") first for any new
# value that would otherwise be re-interpreted as a number, preserving the
# original textual representation (trailing zeros, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.118.10"
$ws.Range("E2").Value = "  -1.56%  "

$ws.Range("D3").Value = "1.895.56"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.18"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5020"
$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3895"
$ws.Range("E8").Value = "  -1.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09233"
$ws.Range("E9").Value = "  -5.69%  "

$ws.Range("E10").Value = "  -2.77%  "

$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.395"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.80"
$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").Value = "1.906.21"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.286"
$ws.Range("E15").Value = "  -3.98%  "

$ws.Range("E17").Value = "  -2.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.31"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.85"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.203"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("D23").Value = "28.187.91"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.320"
$ws.Range("E25").Value = "  +1.89%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.125.90"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.545"
$ws.Range("E27").Value = "  -7.30%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.87"
$ws.Range("E28").Value = "  -2.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.26"
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.81"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.084"
$ws.Range("E31").Value = "  -1.82%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1061"
$ws.Range("E32").Value = "  -1.30%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.605"
$ws.Range("E33").Value = "  -2.31%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.618"
$ws.Range("E34").Value = "  -0.52%  "

$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.571"
$ws.Range("E35").Value = "  -3.00%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.346"
$ws.Range("E36").Value = "  +13.13%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06603"
$ws.Range("E37").Value = "  -3.01%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02397"
$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2206"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.225"
$ws.Range("E40").Value = "  -3.56%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6462"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.974"
$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.39"
$ws.Range("E43").Value = "  -3.08%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6092"
$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.30"
$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.308"
$ws.Range("E47").Value = "  +1.72%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.691"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.006"
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.20"
$ws.Range("E50").Value = "  -2.29%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.200"
$ws.Range("E51").Value = "  -1.11%  "
